# Tier1.xlsx template update:
#  - change large_cnv column layout (insert OMIM_Phenotype_ID column, reorder
#    Gene/Gene-Locus-MIM columns) and hide the Gene/Omim columns used only
#    for lookups (skip gene not found in db)
#  - make large_cnv the active sheet/tab instead of exon_cnv

$wb = $excel.ActiveWorkbook

$wsExon  = $wb.Worksheets.Item("exon_cnv")
$wsLarge = $wb.Worksheets.Item("large_cnv")

# --- large_cnv: update header row (AB1:AH1) ------------------------------
# Before: AB=Gene, AC=OMIM, AD=DiseaseNameEN, AE=DiseaseNameCH, AF=AliasEN,
#         AG=Location, AH=Gene/Locus MIM number
# After:  AB=Gene/Locus MIM number, AC=Gene, AD=OMIM_Phenotype_ID (new),
#         AE=DiseaseNameEN, AF=DiseaseNameCH, AG=AliasEN, AH=Location
$wsLarge.Range("AB1").Value = "Gene/Locus MIM number"
$wsLarge.Range("AC1").Value = "Gene"
$wsLarge.Range("AD1").Value = "OMIM_Phenotype_ID"
$wsLarge.Range("AE1").Value = "DiseaseNameEN"
$wsLarge.Range("AF1").Value = "DiseaseNameCH"
$wsLarge.Range("AG1").Value = "AliasEN"
$wsLarge.Range("AH1").Value = "Location"

# --- large_cnv: hide lookup helper columns L (12) and N (14) ------------
$wsLarge.Columns.Item(12).ColumnWidth = -0.7142857142857143
$wsLarge.Columns.Item(12).Hidden = $true
$wsLarge.Columns.Item(14).ColumnWidth = -0.7142857142857143
$wsLarge.Columns.Item(14).Hidden = $true

# --- switch the active/selected tab from exon_cnv to large_cnv ----------
$wsLarge.Activate()
$excel.ActiveWindow.ScrollColumn = 9
$excel.ActiveWindow.ScrollRow = 1
$wsLarge.Range("AB1").Select() | Out-Null
